# Add TABLE_CONFIG column width specifications to solution briefing tables.
#
# Slide 5 - "Timeline & Milestones" table: resize columns to
#   Phase No 10%, Phase Description 25%, Timeline 15%, Key Deliverables 50%.
# Slide 8 - "Investment Summary" table: resize columns to
#   Cost Category 20%, Year 1 List 10%, Provider/Partner Credits 21%,
#   Year 1 Net 14%, Year 2 11%, Year 3 11%, 3-Year Total 13%.
#
# Column widths are expressed in EMU in the target OOXML; PowerPoint's
# object model works in points, where 1 pt = 12700 EMU. Setting each
# Table.Columns(i).Width also reflows the parent graphicFrame's overall
# Width to the (rounded) sum of the column widths, matching the diff's
# updated <p:xfrm><a:ext cx="..."/> values automatically.

$p = $ppt.ActivePresentation

# --- Slide 5: Timeline & Milestones -----------------------------------
$slideTimeline = $p.Slides.Item(5)
$tblShapeTimeline = $slideTimeline.Shapes.Item(3)
$tblTimeline = $tblShapeTimeline.Table

$timelineWidthsEmu = @(871093, 2177733, 1306639, 4355466)
for ($c = 1; $c -le $tblTimeline.Columns.Count; $c++) {
    $tblTimeline.Columns.Item($c).Width = $timelineWidthsEmu[$c - 1] / 12700
}

# --- Slide 8: Investment Summary ---------------------------------------
$slideInvestment = $p.Slides.Item(8)
$tblShapeInvestment = $slideInvestment.Shapes.Item(3)
$tblInvestment = $tblShapeInvestment.Table

$investmentWidthsEmu = @(1742186, 871093, 1829295, 1219530, 958202, 958202, 1132421)
for ($c = 1; $c -le $tblInvestment.Columns.Count; $c++) {
    $tblInvestment.Columns.Item($c).Width = $investmentWidthsEmu[$c - 1] / 12700
}
